# added 4wk low sales check
$wb = $excel.ActiveWorkbook

# --- "Forecast Comparison" sheet: update Seasonality Index (column L) ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$seasonality = @{
    2  = 1.03
    3  = 0.85
    4  = 1.06
    5  = 0.88
    6  = 1.02
    7  = 1.07
    8  = 1.09
    9  = 1.09
    10 = 1.15
    11 = 0.92
    12 = 0.82
    13 = 0.88
    14 = 1.08
    15 = 0.83
    16 = 0.99
    17 = 0.95
}

foreach ($row in $seasonality.Keys) {
    $wsForecast.Range("L$row").Value = $seasonality[$row]
}

# --- "Summary" sheet: update total forecast figures (keep as text cells) ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "12"
$wsSummary.Range("B9").ClearFormats()

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "6"
$wsSummary.Range("B10").ClearFormats()

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "3"
$wsSummary.Range("B11").ClearFormats()
